$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Update Spectral Type (ST) values for two rows
$ws.Range("B3").Value = "B0Ib"
$ws.Range("B16").Value = "O5fpe"

# Add 2MASS J, H, K magnitude values for rows that were previously blank
$ws.Range("J2").Value = 8.597
$ws.Range("K2").Value = 8.296
$ws.Range("L2").Value = 8.107

$ws.Range("J3").Value = 13.445
$ws.Range("K3").Value = 13.513
$ws.Range("L3").Value = 13.474

$ws.Range("J4").Value = 14.586
$ws.Range("K4").Value = 14.78
$ws.Range("L4").Value = 14.75

$ws.Range("J5").Value = 5.833
$ws.Range("K5").Value = 5.705
$ws.Range("L5").Value = 5.596

$ws.Range("J7").Value = 6.717
$ws.Range("K7").Value = 6.077
$ws.Range("L7").Value = 5.672

$ws.Range("J8").Value = 10.358
$ws.Range("K8").Value = 9.91
$ws.Range("L8").Value = 9.677

$ws.Range("J9").Value = 5.744
$ws.Range("K9").Value = 5.639
$ws.Range("L9").Value = 5.496

$ws.Range("J11").Value = 13.695
$ws.Range("K11").Value = 13.537
$ws.Range("L11").Value = 13.293

$ws.Range("J12").Value = 6.872
$ws.Range("K12").Value = 6.652
$ws.Range("L12").Value = 6.501

# Update the active selection to match the recorded view state
$ws.Range("J13").Select()
